$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.801.28'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '2.192.53'
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '291.04'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('D6').Value = '86.57'
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('E7').Value = '  -1.59%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '0.462'
$ws.Range('E9').Value = '  -2.07%  '
$ws.Range('D10').Value = '30.00'
$ws.Range('E10').Value = '  -3.57%  '
$ws.Range('D11').Value = '49.95'
$ws.Range('E11').Value = '  +6.12%  '
$ws.Range('E12').Value = '  -1.80%  '
$ws.Range('E13').Value = '  +2.54%  '
$ws.Range('D14').Value = '6.39'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').Value = '2.530.84'
$ws.Range('E15').Value = '  -1.84%  '
$ws.Range('E16').Value = '  -3.07%  '
$ws.Range('D17').Value = '2.192.46'
$ws.Range('E17').Value = '  -1.51%  '
$ws.Range('D18').Value = '0.723'
$ws.Range('E18').Value = '  -1.05%  '
$ws.Range('D19').Value = '39.697.24'
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('E20').Value = '  -1.08%  '
$ws.Range('D21').Value = '11.14'
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('E22').Value = '  -2.01%  '
$ws.Range('D23').Value = '65.07'
$ws.Range('E23').Value = '  -0.62%  '
$ws.Range('D24').Value = '236.19'
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E26').Value = '  -1.28%  '
$ws.Range('E27').Value = '  -2.81%  '
$ws.Range('D28').Value = '23.07'
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('D29').Value = '2.04'
$ws.Range('E29').Value = '  -8.00%  '
$ws.Range('D30').Value = '9.13'
$ws.Range('E30').Value = '  -2.03%  '
$ws.Range('D31').Value = '155.61'
$ws.Range('E31').Value = '  +2.40%  '
$ws.Range('D32').Value = '31.08'
$ws.Range('E32').Value = '  -6.24%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('E34').Value = '  -0.70%  '
$ws.Range('E35').Value = '  -2.42%  '
$ws.Range('E36').Value = '  -2.40%  '
$ws.Range('D37').Value = '2.81'
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = '0.0970'
$ws.Range('E39').Value = '  -2.96%  '
$ws.Range('D40').Value = '15.07'
$ws.Range('E40').Value = '  -7.30%  '
$ws.Range('E41').Value = '  -3.21%  '
$ws.Range('D42').Value = '2.123.90'
$ws.Range('E42').Value = '  +2.68%  '
$ws.Range('E43').Value = '  -2.54%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').Value = '2.10'
$ws.Range('E44').Value = '  -1.22%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '0.0266'
$ws.Range('E45').Value = '  -1.22%  '
$ws.Range('E46').Value = '  -2.32%  '
$ws.Range('D47').Value = '17.05'
$ws.Range('E47').Value = '  -6.50%  '
$ws.Range('E48').Value = '  +2.44%  '
$ws.Range('D49').Value = '2.402.52'
$ws.Range('E49').Value = '  -1.49%  '
$ws.Range('E50').Value = '  +1.17%  '
$ws.Range('D51').Value = '87.80'
$ws.Range('E51').Value = '  -1.58%  '
